$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("data")
$wsCustomer = $wb.Worksheets.Item("CustomerData")

# Add Password column to "data" sheet (sheet2)
$wsData.Range("I1").Value = "Password"
$wsData.Range("I2").Value = 123
$wsData.Range("I3").Value = 1234
$wsData.Range("I4").Value = 12345
$wsData.Range("I5").Value = 123456
$wsData.Range("I6").Value = 12345
$wsData.Range("I7").Value = 1234
$wsData.Range("I8").Value = 123
$wsData.Range("I9").Value = 12
$wsData.Range("I10").Value = 1

# Add Password column to "CustomerData" sheet (sheet3)
$wsCustomer.Range("F1").Value = "Password"
$wsCustomer.Range("F2").Value = 123
$wsCustomer.Range("F3").Value = 1234
$wsCustomer.Range("F4").Value = 12345
$wsCustomer.Range("F5").Value = 123456
$wsCustomer.Range("F6").Value = 12345
$wsCustomer.Range("F7").Value = 1234
$wsCustomer.Range("F8").Value = 123
$wsCustomer.Range("F9").Value = 12
$wsCustomer.Range("F10").Value = 1
$wsCustomer.Range("F11").Value = 1

# Update selections / active sheet state
$wsCustomer.Activate()
$wsCustomer.Range("F1:F11").Select()

$wsData.Activate()
$wsData.Range("G13").Select()
